$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format the price/volume columns as Text so numeric-looking strings
# (e.g. "341.49", "1.001", "18.00") are preserved verbatim instead of being
# auto-converted to floating point numbers by the smart-entry heuristics.
$priceRange = $ws.Range("D2:E51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = '29.796.10'
$ws.Range("E2").Value = '  +7.15%  '
$ws.Range("D3").Value = '1.950.01'
$ws.Range("E3").Value = '  +5.39%  '
$ws.Range("E4").Value = '  -0.57%  '
$ws.Range("D5").Value = '341.49'
$ws.Range("E5").Value = '  +1.99%  '
$ws.Range("E6").Value = '  -0.43%  '
$ws.Range("D7").Value = '0.4791'
$ws.Range("E7").Value = '  +3.01%  '
$ws.Range("D8").Value = '0.4134'
$ws.Range("E8").Value = '  +6.84%  '
$ws.Range("D9").Value = '47.84'
$ws.Range("E9").Value = '  +2.03%  '
$ws.Range("D10").Value = '0.08235'
$ws.Range("E10").Value = '  +4.10%  '
$ws.Range("D11").Value = '1.038'
$ws.Range("E11").Value = '  +7.15%  '
$ws.Range("D12").Value = '22.78'
$ws.Range("E12").Value = '  +6.68%  '
$ws.Range("D13").Value = '1.964.75'
$ws.Range("E13").Value = '  +6.59%  '
$ws.Range("D14").Value = '6.156'
$ws.Range("E14").Value = '  +4.41%  '
$ws.Range("E15").Value = '  +3.22%  '
$ws.Range("D16").Value = '91.92'
$ws.Range("E16").Value = '  +1.89%  '
$ws.Range("D17").Value = '1.001'
$ws.Range("E17").Value = '  -0.54%  '
$ws.Range("D18").Value = '0.00001059'
$ws.Range("E18").Value = '  +2.94%  '
$ws.Range("D19").Value = '0.06673'
$ws.Range("E19").Value = '  +0.98%  '
$ws.Range("D20").Value = '18.00'
$ws.Range("E20").Value = '  +3.68%  '
$ws.Range("E21").Value = '  -0.54%  '
$ws.Range("D22").Value = '29.765.61'
$ws.Range("E22").Value = '  +7.04%  '
$ws.Range("D23").Value = '5.588'
$ws.Range("E23").Value = '  +4.46%  '
$ws.Range("D24").Value = '11.27'
$ws.Range("E24").Value = '  +3.83%  '
$ws.Range("D25").Value = '2.288'
$ws.Range("E25").Value = '  -0.32%  '
$ws.Range("D26").Value = '2.183.17'
$ws.Range("E26").Value = '  +5.70%  '
$ws.Range("D27").Value = '161.15'
$ws.Range("E27").Value = '  +1.66%  '
$ws.Range("D28").Value = '20.24'
$ws.Range("E28").Value = '  +3.94%  '
$ws.Range("E29").Value = '  +5.12%  '
$ws.Range("D30").Value = '5.660'
$ws.Range("E30").Value = '  +5.42%  '
$ws.Range("D31").Value = '122.94'
$ws.Range("E31").Value = '  +3.47%  '
$ws.Range("D32").Value = '1.008'
$ws.Range("E32").Value = '  +6.42%  '
$ws.Range("D33").Value = '0.09663'
$ws.Range("E33").Value = '  +2.60%  '
$ws.Range("E34").Value = '  +11.14%  '
$ws.Range("D35").Value = '3.685'
$ws.Range("E35").Value = '  +2.66%  '
$ws.Range("D36").Value = '5.502'
$ws.Range("E36").Value = '  +4.53%  '
$ws.Range("D37").Value = '0.06259'
$ws.Range("E37").Value = '  +3.92%  '
$ws.Range("D38").Value = '0.02316'
$ws.Range("E38").Value = '  +4.42%  '
$ws.Range("D39").Value = '8.504'
$ws.Range("E39").Value = '  +2.89%  '
$ws.Range("D40").Value = '1.186'
$ws.Range("E40").Value = '  +2.63%  '
$ws.Range("D41").Value = '0.6083'
$ws.Range("E41").Value = '  +4.62%  '
$ws.Range("D42").Value = '10.73'
$ws.Range("E42").Value = '  +6.53%  '
$ws.Range("D43").Value = '0.1895'
$ws.Range("E43").Value = '  +2.85%  '
$ws.Range("D44").Value = '1.000'
$ws.Range("E44").Value = '  -0.42%  '
$ws.Range("D45").Value = '2.396'
$ws.Range("E45").Value = '  +32.44%  '
$ws.Range("E46").Value = '  -0.38%  '
$ws.Range("D47").Value = '0.5704'
$ws.Range("E47").Value = '  +4.63%  '
$ws.Range("D48").Value = '12.51'
$ws.Range("E48").Value = '  +4.83%  '
$ws.Range("D49").Value = '0.07421'
$ws.Range("E49").Value = '  +8.33%  '
$ws.Range("E50").Value = '  +3.26%  '
$ws.Range("D51").Value = '113.09'
$ws.Range("E51").Value = '  +1.96%  '

# Restore the original (default/no explicit) cell formatting now that the
# text values are committed, so styling matches the untouched rows/cells.
$priceRange.ClearFormats()

Write-Host "Updated 92 cells"
